$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (E3/F3): update the scenario/expected-result text
$ws.Range("E3").Value = "Verify  OTP should expire after cofigured time"
$ws.Range("F3").Value = "OTP should get expire"

# Move the saved cursor/selection to H5 (matches the author's re-save state)
$ws.Range("H5").Select()
